# Apply "Preparar clase Tag para la entrega" edits to controlDedicacion.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26: dedicated time changes from 45 to 60 minutes ---
$ws.Range("C26").Value = 60

# --- Row 27: new entry (15 min, Documentación, "Tarea 12") ---
$ws.Range("C27").Value = 15
$ws.Range("G27").Value = "Documentación"
$ws.Range("H27").Value = "Tarea 12"

# --- Row 28: new entry (15 min, Otros, "Preparar archivos para entrega") ---
$ws.Range("C28").Value = 15
$ws.Range("G28").Value = "Otros"
$ws.Range("H28").Value = "Preparar archivos para entrega"

# --- Row 29: new entry (15 min, Pruebas, "Tareas 11 y 12") ---
$ws.Range("C29").Value = 15
$ws.Range("G29").Value = "Pruebas"
$ws.Range("H29").Value = "Tareas 11 y 12"

# --- Update the view: scroll position and active selection moved to H29 ---
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H29").Select()
